# Applies the "amcr/externi.docx" style-sheet update:
#   1. Adds a new "Abstract Title" paragraph style (based on Normal,
#      followed by Abstract).
#   2. Tightens the existing "Abstract" style's space-before from 15pt
#      (300 twips) to 5pt (100 twips).
#   3. Adds a new "Footnote Block Text" paragraph style (based on
#      Footnote Text, followed by Footnote Text).

$d = $word.ActiveDocument

# --- 1. New "Abstract Title" style ------------------------------------
$abstractTitle = $d.Styles.Add("Abstract Title", 1)
$abstractTitle.BaseStyle = "Normal"
$abstractTitle.NextParagraphStyle = "Abstract"
$abstractTitle.QuickStyle = $true

$atPf = $abstractTitle.ParagraphFormat
$atPf.KeepWithNext = $true
$atPf.KeepTogether = $true
$atPf.Alignment = 1        # wdAlignParagraphCenter
$atPf.SpaceBefore = 15     # 300 twips
$atPf.SpaceAfter = 0

$atFont = $abstractTitle.Font
$atFont.Size = 10          # sz = 20 (half-points)
$atFont.SizeBi = 10        # szCs = 20
$atFont.Bold = $true
$atFont.Color = 9067060    # RGB(0x34,0x5A,0x8A) -> OLE BGR

# --- 2. Tighten "Abstract" style's space-before ------------------------
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5   # 100 twips

# --- 3. New "Footnote Block Text" style --------------------------------
# (Styles.Add always mints a w:customStyle="1" entry on this object model,
# same as genuine Word COM automation for any non-reserved style id.)
$footnoteBlockText = $d.Styles.Add("Footnote Block Text", 1)
$footnoteBlockText.BaseStyle = "Footnote Text"
$footnoteBlockText.NextParagraphStyle = "Footnote Text"
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true
$footnoteBlockText.QuickStyle = $true

$fbtPf = $footnoteBlockText.ParagraphFormat
$fbtPf.SpaceBefore = 5     # 100 twips
$fbtPf.SpaceAfter = 5      # 100 twips
$fbtPf.FirstLineIndent = 0
$fbtPf.LeftIndent = 24     # 480 twips
$fbtPf.RightIndent = 24    # 480 twips

Write-Output "Styles updated: Abstract Title, Abstract, Footnote Block Text"
